$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string labelled cells (order matters for sharedStrings.xml layout)
$ws.Range("D2").Value = "One 30-minute download will take (secs)"
$ws.Range("D3").Value = 11

# Update J3 download duration in seconds
$ws.Range("J3").Value = 1161

# Headers for the secs/hours/days summary block, row 6 (Secs added before Criteria block)
$ws.Range("M6").Value = "Secs"
$ws.Range("N6").Value = "Hours"
$ws.Range("O6").Value = "Days"
$ws.Range("O6").Font.Bold = $true

# Criteria block
$ws.Range("J4").Value = "Criteria:"
$ws.Range("J5").Value = "20 day average volume: >400,000"
$ws.Range("J6").Value = "IPO: > 1 year ago"
$ws.Range("J7").Value = "Exchanges: all US - OTC"
$ws.Range("J8").Value = "Price range: 1 - 499"
$ws.Range("J9").Value = "Securty Type: != MLP Ltd Part"

# Row 7 computed values
$ws.Range("M7").Formula = "=M3*D3"
$ws.Range("N7").Formula = "=M7/3600"
$ws.Range("O7").Formula = "=N7/24"
$ws.Range("O7").Font.Bold = $true

# Column D width (closest reachable value given engine's pixel rounding)
$ws.Range("D1").ColumnWidth = 32.28515625

# Selection
$ws.Range("J3").Select() | Out-Null

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
